$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (incl. date number format) from the last existing
# data row (row 7) onto the new row 8 before populating values, so the
# new row's style indexes line up with the existing styles (no new
# numFmt/cellXf entries get minted).
$ws.Range("A7:H7").Copy($ws.Range("A8:H8"))

$ws.Range("A8").Value = 9309.17
$ws.Range("B8").Value = 9407
$ws.Range("C8").Value = 109.08
$ws.Range("D8").Value = 107.95
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -1.04
$ws.Range("G8").Value = 42612.672951388886
$ws.Range("H8").Value = $false
